$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '68.500.63'
Set-TextValue $ws.Range("E2") '  +1.35%  '

Set-TextValue $ws.Range("D3") '3.921.48'
Set-TextValue $ws.Range("E3") '  +1.19%  '

Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.14%  '

Set-TextValue $ws.Range("D5") '484.04'
Set-TextValue $ws.Range("E5") '  +4.75%  '

Set-TextValue $ws.Range("D6") '147.72'
Set-TextValue $ws.Range("E6") '  -0.38%  '

Set-TextValue $ws.Range("E7") '  -1.38%  '

Set-TextValue $ws.Range("E8") '  -0.01%  '

Set-TextValue $ws.Range("E9") '  -4.25%  '

Set-TextValue $ws.Range("E10") '  +7.60%  '

Set-TextValue $ws.Range("D11") '0.0000355'
Set-TextValue $ws.Range("E11") '  +12.05%  '

Set-TextValue $ws.Range("D12") '42.40'
Set-TextValue $ws.Range("E12") '  -4.02%  '

Set-TextValue $ws.Range("D13") '10.49'
Set-TextValue $ws.Range("E13") '  +0.22%  '

Set-TextValue $ws.Range("D14") '4.547.29'
Set-TextValue $ws.Range("E14") '  +1.39%  '

Set-TextValue $ws.Range("B15") 'Uniswap'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D15") '14.61'
Set-TextValue $ws.Range("E15") '  -1.72%  '

Set-TextValue $ws.Range("B16") 'WrappedEther'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D16") '3.958.86'
Set-TextValue $ws.Range("E16") '  +2.09%  '

Set-TextValue $ws.Range("E17") '  -0.46%  '

Set-TextValue $ws.Range("E18") '  -2.37%  '

Set-TextValue $ws.Range("D19") '1.13'
Set-TextValue $ws.Range("E19") '  -3.25%  '

Set-TextValue $ws.Range("D20") '68.660.07'
Set-TextValue $ws.Range("E20") '  +1.35%  '

Set-TextValue $ws.Range("D21") '430.85'
Set-TextValue $ws.Range("E21") '  +0.15%  '

Set-TextValue $ws.Range("B22") 'ImmutableX'
Set-TextValue $ws.Range("C22") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D22") '3.33'
Set-TextValue $ws.Range("E22") '  +1.15%  '

Set-TextValue $ws.Range("B23") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C23") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D23") '14.49'
Set-TextValue $ws.Range("E23") '  -3.01%  '

Set-TextValue $ws.Range("E24") '  -1.06%  '

Set-TextValue $ws.Range("D25") '11.35'
Set-TextValue $ws.Range("E25") '  +11.38%  '

Set-TextValue $ws.Range("D26") '3.57'
Set-TextValue $ws.Range("E26") '  +0.50%  '

Set-TextValue $ws.Range("D27") '10.55'
Set-TextValue $ws.Range("E27") '  +1.46%  '

Set-TextValue $ws.Range("D28") '38.11'
Set-TextValue $ws.Range("E28") '  +0.65%  '

Set-TextValue $ws.Range("D29") '5.89'
Set-TextValue $ws.Range("E29") '  +7.23%  '

Set-TextValue $ws.Range("D30") '705.08'
Set-TextValue $ws.Range("E30") '  -6.27%  '

Set-TextValue $ws.Range("D31") '13.21'
Set-TextValue $ws.Range("E31") '  -3.91%  '

Set-TextValue $ws.Range("E32") '  -4.85%  '

Set-TextValue $ws.Range("E33") '  +2.84%  '

Set-TextValue $ws.Range("D34") '0.0₃0898'
Set-TextValue $ws.Range("E34") '  +31.65%  '

Set-TextValue $ws.Range("D35") '41.34'
Set-TextValue $ws.Range("E35") '  -5.41%  '

Set-TextValue $ws.Range("D36") '58.44'
Set-TextValue $ws.Range("E36") '  +1.49%  '

Set-TextValue $ws.Range("D37") '0.150'
Set-TextValue $ws.Range("E37") '  -8.17%  '

Set-TextValue $ws.Range("B38") 'Dai'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D38") '0.999'
Set-TextValue $ws.Range("E38") '  -0.14%  '

Set-TextValue $ws.Range("B39") 'NEARProtocol'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D39") '5.45'
Set-TextValue $ws.Range("E39") '  -1.95%  '

Set-TextValue $ws.Range("D40") '2.87'
Set-TextValue $ws.Range("E40") '  +8.87%  '

Set-TextValue $ws.Range("E41") '  -2.58%  '

Set-TextValue $ws.Range("D42") '3.07'
Set-TextValue $ws.Range("E42") '  +12.20%  '

Set-TextValue $ws.Range("E43") '  +1.09%  '

Set-TextValue $ws.Range("D44") '0.343'
Set-TextValue $ws.Range("E44") '  -2.89%  '

Set-TextValue $ws.Range("B45") 'FirstDigitalUSD'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D45") '1.00'
Set-TextValue $ws.Range("E45") '  -0.01%  '

Set-TextValue $ws.Range("B46") 'Stellar'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D46") '0.140'
Set-TextValue $ws.Range("E46") '  -0.72%  '

Set-TextValue $ws.Range("D47") '3.25'
Set-TextValue $ws.Range("E47") '  -1.17%  '

Set-TextValue $ws.Range("B48") 'LidoDAOToken'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D48") '3.40'
Set-TextValue $ws.Range("E48") '  -1.85%  '

Set-TextValue $ws.Range("B49") 'ARBITRUM'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D49") '2.16'
Set-TextValue $ws.Range("E49") '  +0.74%  '

Set-TextValue $ws.Range("D50") '147.06'
Set-TextValue $ws.Range("E50") '  +1.24%  '

Set-TextValue $ws.Range("D51") '2.83'
Set-TextValue $ws.Range("E51") '  -2.50%  '
